$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: volume/issue number and report week date range ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- Weekly crime statistics table (rows 15-31) ---
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -78.787878787878
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 136
$ws.Range("K16").Value = -29.411764705882
$ws.Range("L16").Value = -30.434782608695
$ws.Range("M16").Value = -48.387096774193
$ws.Range("N16").Value = -85.756676557863
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -10.526315789473
$ws.Range("I17").Value = 198
$ws.Range("J17").Value = 240
$ws.Range("K17").Value = -17.5
$ws.Range("L17").Value = -1
$ws.Range("M17").Value = 32.885906040268
$ws.Range("N17").Value = -66.153846153846
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 61
$ws.Range("J18").Value = 74
$ws.Range("K18").Value = -17.567567567567
$ws.Range("L18").Value = -32.222222222222
$ws.Range("M18").Value = -30.681818181818
$ws.Range("N18").Value = -90.976331360946
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 55.555555555555
$ws.Range("F19").Value = 43
$ws.Range("H19").Value = 13.157894736842
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 402
$ws.Range("K19").Value = -19.154228855721
$ws.Range("L19").Value = 7.97342192691
$ws.Range("M19").Value = 35.983263598326
$ws.Range("N19").Value = -12.398921832884
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 2
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 24
$ws.Range("K20").Value = -48.936170212766
$ws.Range("L20").Value = -57.142857142857
$ws.Range("M20").Value = 41.176470588235
$ws.Range("N20").Value = -78.947368421052
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -3.846153846153
$ws.Range("I21").Value = 716
$ws.Range("J21").Value = 911
$ws.Range("K21").Value = -21.405049396267
$ws.Range("L21").Value = -9.937106918238
$ws.Range("M21").Value = 3.318903318903
$ws.Range("N21").Value = -71.152296535052
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -17.647058823529
$ws.Range("M22").Value = 40
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 66
$ws.Range("K23").Value = -28.787878787878
$ws.Range("L23").Value = -33.802816901408
$ws.Range("M23").Value = 34.285714285714
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 91.666666666666
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 11.235955056179
$ws.Range("I24").Value = 1189
$ws.Range("J24").Value = 976
$ws.Range("K24").Value = 21.823770491803
$ws.Range("L24").Value = 19.37751004016
$ws.Range("M24").Value = 54.817708333333
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 760
$ws.Range("J25").Value = 594
$ws.Range("K25").Value = 27.946127946127
$ws.Range("L25").Value = 31.03448275862
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = -34.146341463414
$ws.Range("I26").Value = 335
$ws.Range("J26").Value = 350
$ws.Range("K26").Value = -4.285714285714
$ws.Range("L26").Value = 4.36137071651
$ws.Range("M26").Value = -17.690417690417
$ws.Range("L27").Value = -30.76923076923
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 3
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = -18.421052631578
$ws.Range("N29").Value = -91.208791208791
$ws.Range("N30").Value = -91.25
$ws.Range("D31").Copy($ws.Range("C31"))
